$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2018
$ws.Range("B2").Value = 3.269218759757697
$ws.Range("A3").Value = 2019
$ws.Range("B3").Value = 4.190783698126683
$ws.Range("A4").Value = 2020
$ws.Range("B4").Value = 5.176704112308157
$ws.Range("A5").Value = 2021
$ws.Range("B5").Value = 3.862804326340846
$ws.Range("A6:B6").Delete()

$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$s = $chart.SeriesCollection(1)
$s.Formula = '=SERIES(Sheet1!$B$1,Sheet1!$A$2:$A$5,Sheet1!$B$2:$B$5,1)'
$wb.RefreshAll()
Write-Host $s.Formula
